# SMPTE test pattern workbook: add a "Colors" legend (named-rectangle color
# swatches) in columns G:H of Sheet1, plus the supporting shared strings,
# a new left-aligned cell style, and the updated view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column G (width, in "characters") + left alignment style for the
#     whole color-legend block (skip row 18, which stays empty/unused). ---
$ws.Columns("G").ColumnWidth = 20
$ws.Range("G2:G17").HorizontalAlignment = -4131   # xlHAlignLeft
$ws.Range("G19:G30").HorizontalAlignment = -4131  # xlHAlignLeft
$ws.Range("H17").HorizontalAlignment = -4131      # xlHAlignLeft

# --- Header (default style) ---
$ws.Range("G1").Value = "Colors"

# --- Color legend rows 2-17 (named rectangles + sampled pixel colors) ---
$ws.Range("G2").Value = "0xff696969"
$ws.Range("G3").Value = "c1c1c1"
$ws.Range("G4").Value = "c1c100"
$ws.Range("G5").Value = "00c1c1"
$ws.Range("G6").Value = "00c100"
$ws.Range("G7").Value = "c100c1"
$ws.Range("G8").Value = "c10000"
$ws.Range("G9").Value = "0000c1"
$ws.Range("G10").Value = 696969
$ws.Range("G11").Value = "00ffff"
$ws.Range("G12").Value = 52550
$ws.Range("G13").Value = "c1c1c1"
$ws.Range("G14").Value = "0000ff"
$ws.Range("G15").Value = "ffff00"
$ws.Range("G16").Value = "36056d"
$ws.Range("G17").Value = 50505
$ws.Range("H17").Value = "fdfdfd"

# --- Rows 19-30 (row 18 intentionally left blank) ---
$ws.Range("G19").Value = "ff0000"
$ws.Range("G20").Value = "2b2b2b"
$ws.Range("G21").Value = 50505
$ws.Range("G22").Value = "ffffff"
$ws.Range("G23").Value = 50505
$ws.Range("G24").Value = 0
$ws.Range("G25").Value = 50505
$ws.Range("G26").Value = "0a0a0a"
$ws.Range("G27").Value = 50505
$ws.Range("G28").Value = "0d0d0d"
$ws.Range("G29").Value = 50505
$ws.Range("G30").Value = "2b2b2b"

# --- View state: scroll so the legend is visible, select H17 ---
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("H17").Select()
